# Auto-generated edits applying the diff to Alpha_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2128.5715
$ws.Range("I70").Value = 2140
$ws.Range("K70").Value = 6420
$ws.Range("M70").Value = -6150

$ws.Range("H73").Value = 2128.5715
$ws.Range("I73").Value = 2140
$ws.Range("K73").Value = 6420
$ws.Range("M73").Value = -5484

$ws.Range("H132").Value = 57321.332
$ws.Range("I132").Value = 60105.234
$ws.Range("J132").Value = 9995
$ws.Range("K132").Value = 180315.702
$ws.Range("L132").Value = 29985
$ws.Range("M132").Value = -177785.702
$ws.Range("N132").Value = -35045

$ws.Range("H137").Value = 988.26666
$ws.Range("I137").Value = 979
$ws.Range("J137").Value = 1025.3334
$ws.Range("K137").Value = 2937
$ws.Range("L137").Value = 3076.0002
$ws.Range("M137").Value = -387
$ws.Range("N137").Value = -8176.0002

$ws.Range("H138").Value = 7945.933
$ws.Range("I138").Value = 5027.2
$ws.Range("J138").Value = 8779.857
$ws.Range("K138").Value = 15081.6
$ws.Range("L138").Value = 26339.571
$ws.Range("M138").Value = -9941.599999999999
$ws.Range("N138").Value = -36619.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 18001
$ws.Range("I12").Value = 4000
$ws.Range("J12").Value = 60004
$ws.Range("K12").Value = 4000
$ws.Range("L12").Value = 60004
$ws.Range("M12").Value = -3827
$ws.Range("N12").Value = -60350

$ws.Range("H45").Value = 2572.0908
$ws.Range("I45").Value = 2175.5
$ws.Range("J45").Value = 3629.6667
$ws.Range("K45").Value = 2175.5
$ws.Range("L45").Value = 3629.6667
$ws.Range("M45").Value = -1798.5
$ws.Range("N45").Value = -4383.6667

$ws.Range("H74").Value = 490013.5
$ws.Range("I74").Value = 2912.3408
$ws.Range("K74").Value = 2912.3408
$ws.Range("M74").Value = -2038.3408

$ws.Range("H77").Value = 490013.5
$ws.Range("I77").Value = 2912.3408
$ws.Range("K77").Value = 14561.704
$ws.Range("M77").Value = -10193.704

$ws.Range("H110").Value = 10404
$ws.Range("I110").Value = 9998.333000000001
$ws.Range("K110").Value = 9998.333000000001
$ws.Range("M110").Value = -7953.333000000001

$ws.Range("H122").Value = 5013.7617
$ws.Range("I122").Value = 4583
$ws.Range("K122").Value = 13749
$ws.Range("M122").Value = -11299

$ws.Range("H125").Value = 83722
$ws.Range("J125").Value = 83722
$ws.Range("L125").Value = 83722
$ws.Range("N125").Value = -93562

$ws.Range("H132").Value = 20837600
$ws.Range("I132").Value = 3517.9092
$ws.Range("K132").Value = 10553.7276
$ws.Range("M132").Value = -8023.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 23082.2
$ws.Range("J94").Value = 2677.25
$ws.Range("L94").Value = 2677.25
$ws.Range("N94").Value = -3579.25

$ws.Range("H99").Value = 1593.5
$ws.Range("I99").Value = 1593.5
$ws.Range("K99").Value = 1593.5
$ws.Range("M99").Value = -95.5

$ws.Range("H105").Value = 6148
$ws.Range("I105").Value = 6148
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 6148
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -4401
$ws.Range("N105").Value = ""

$ws.Range("H107").Value = 14021.895
$ws.Range("I107").Value = 1009.7692
$ws.Range("K107").Value = 1009.7692
$ws.Range("M107").Value = 910.2308

$ws.Range("H131").Value = 300000
$ws.Range("J131").Value = 300000
$ws.Range("L131").Value = 300000
$ws.Range("N131").Value = -310080

$ws.Range("H134").Value = 19232104
$ws.Range("I134").Value = 20834644
$ws.Range("K134").Value = 62503932
$ws.Range("M134").Value = -62501397

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3122.7144
$ws.Range("I16").Value = 3143.3333
$ws.Range("J16").Value = 2999
$ws.Range("K16").Value = 3143.3333
$ws.Range("L16").Value = 2999
$ws.Range("M16").Value = -2856.3333
$ws.Range("N16").Value = -3573

$ws.Range("H31").Value = 2182.853
$ws.Range("J31").Value = 1936.6522
$ws.Range("L31").Value = 1936.6522
$ws.Range("N31").Value = -2526.6522

$ws.Range("H34").Value = 2182.853
$ws.Range("J34").Value = 1936.6522
$ws.Range("L34").Value = 1936.6522
$ws.Range("N34").Value = -2340.6522

$ws.Range("H113").Value = 3122.7144
$ws.Range("I113").Value = 3143.3333
$ws.Range("J113").Value = 2999
$ws.Range("K113").Value = 3143.3333
$ws.Range("L113").Value = 2999
$ws.Range("M113").Value = -973.3332999999998
$ws.Range("N113").Value = -7339

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 969.8333
$ws.Range("J12").Value = 1667
$ws.Range("L12").Value = 5001
$ws.Range("N12").Value = -5347

$ws.Range("H68").Value = 1528.1364
$ws.Range("J68").Value = 1567.8889
$ws.Range("L68").Value = 4703.6667
$ws.Range("N68").Value = -6325.6667

$ws.Range("H71").Value = 1528.1364
$ws.Range("J71").Value = 1567.8889
$ws.Range("L71").Value = 14111.0001
$ws.Range("N71").Value = -22223.0001

$ws.Range("H121").Value = 95856.17999999999
$ws.Range("J121").Value = 5935.4443
$ws.Range("L121").Value = 17806.3329
$ws.Range("N121").Value = -20426.3329

$ws.Range("H134").Value = 10456.8
$ws.Range("I134").Value = 10456.8
$ws.Range("K134").Value = 31370.4
$ws.Range("M134").Value = -26300.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3007.45
$ws.Range("I102").Value = 3241.8333
$ws.Range("K102").Value = 3241.8333
$ws.Range("M102").Value = -1619.8333

$ws.Range("H122").Value = 3908.2856
$ws.Range("I122").Value = 4202.3
$ws.Range("J122").Value = 3173.25
$ws.Range("K122").Value = 12606.9
$ws.Range("L122").Value = 9519.75
$ws.Range("M122").Value = -10156.9
$ws.Range("N122").Value = -14419.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4789.4
$ws.Range("I22").Value = 817.1667
$ws.Range("K22").Value = 817.1667
$ws.Range("M22").Value = -522.1667

$ws.Range("H27").Value = 4789.4
$ws.Range("I27").Value = 817.1667
$ws.Range("K27").Value = 817.1667
$ws.Range("M27").Value = -710.1667

$ws.Range("H46").Value = 2750.3333
$ws.Range("J46").Value = 2750.3333
$ws.Range("L46").Value = 2750.3333
$ws.Range("N46").Value = -3126.3333

$ws.Range("H61").Value = 2513.3635
$ws.Range("I61").Value = 2564.7
$ws.Range("K61").Value = 2564.7
$ws.Range("M61").Value = -2362.7

$ws.Range("H113").Value = 2513.3635
$ws.Range("I113").Value = 2564.7
$ws.Range("K113").Value = 2564.7
$ws.Range("M113").Value = -394.6999999999998

$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999

$ws.Range("H132").Value = 2250
$ws.Range("I132").Value = 2250
$ws.Range("K132").Value = 6750
$ws.Range("M132").Value = -4220

$ws.Range("H136").Value = 100003656
$ws.Range("I136").Value = 2914
$ws.Range("K136").Value = 8742
$ws.Range("M136").Value = -6192

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 10023.333
$ws.Range("I51").Value = 10023.333
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 10023.333
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -9513.333000000001
$ws.Range("N51").Value = ""

$ws.Range("H62").Value = 3500
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3500
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740
